$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.893.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "'3.438.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D5").Value = "'573.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "'175.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.52%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  +4.70%  "
$ws.Range("D11").Value = "'55.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'0.0000273"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "'9.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").Value = "'3.982.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "'3.447.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "'18.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'11.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'64.885.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.03%  "
$ws.Range("D20").Value = "'0.991"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "'408.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.52%  "
$ws.Range("D22").Value = "'4.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").Value = "'4.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("D24").Value = "'83.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'13.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.85%  "
$ws.Range("D26").Value = "'10.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").Value = "'2.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("D28").Value = "'6.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").Value = "'8.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "'29.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "'6.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "'587.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.68%  "
$ws.Range("D33").Value = "'11.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("D35").Value = "'59.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").Value = "'0.0₃0772"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'36.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.01%  "
$ws.Range("D41").Value = "'0.378"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").Value = "'3.182.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.69%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'2.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.63%  "
$ws.Range("D47").Value = "'0.0410"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'0.131"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  -4.36%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'8.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'136.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.04%  "
